$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the capitalization of the misplaced entry (it was already alphabetized
# at row 9 in this sheet) and flag it, along with the other two out-of-place
# names at the bottom of the list, in red.
$ws.Range("A9").Value2 = "Jhonatan Giraldo"

$redColor = 255
$ws.Range("A9").Font.Color = $redColor
$ws.Range("A28:A30").Font.Color = $redColor

# Remove the two trailing empty rows.
$ws.Range("A31:A32").EntireRow.Delete()

# Update the active cell selection to match the saved workbook state.
$ws.Range("C10").Select()
